# Add two new Mac-Addresses: 10 new device_master rows (two full cycles of
# the 5 device types: Finger Print Scanner, IRIS Scanner, Web Camera,
# Document Scanner, Printer), continuing the existing table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# columns: id, name, mac_address, serial_num, dspec_id
$newRows = @(
    @(3000166, "Finger Print Scanner 30", "D6-15-AC-80-6B-86", "BS563Q2230814", 165),
    @(3000167, "IRIS Scanner 30",         "6D-58-E2-DF-74-34", "BS563Q2230815", 327),
    @(3000168, "Web Camera 30",           "E2-A8-56-86-15-30", "BS563Q2230816", 736),
    @(3000169, "Document Scanner 30",     "72-E8-B9-FD-63-65", "BS563Q2230817", 801),
    @(3000170, "Printer 30",              "D3-F3-A4-50-AD-12", "BS563Q2230818", 920),
    @(3000171, "Finger Print Scanner 31", "06-16-D0-0B-A6-E4", "BS563Q2230819", 165),
    @(3000172, "IRIS Scanner 31",         "21-78-45-AC-E9-20", "BS563Q2230820", 327),
    @(3000173, "Web Camera 31",           "3C-E8-87-99-DB-FA", "BS563Q2230821", 736),
    @(3000174, "Document Scanner 31",     "BF-55-53-98-40-08", "BS563Q2230822", 801),
    @(3000175, "Printer 31",              "5A-43-36-46-22-EB", "BS563Q2230823", 920)
)

$startRow = 147
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value2 = $data[0]   # id
    $ws.Cells.Item($r, 2).Value2 = $data[1]   # name
    $ws.Cells.Item($r, 3).Value2 = $data[2]   # mac_address
    $ws.Cells.Item($r, 4).Value2 = $data[3]   # serial_num
    $ws.Cells.Item($r, 6).Value2 = $data[4]   # dspec_id
    $ws.Cells.Item($r, 7).Value2 = "eng"      # lang_code
    $ws.Cells.Item($r, 8).Value2 = $true      # is_active
    $ws.Cells.Item($r, 8).HorizontalAlignment = -4131   # xlLeft, matches other rows' style
    $ws.Cells.Item($r, 9).Value2 = "superadmin"  # cr_by
    $ws.Cells.Item($r, 10).Value2 = "now()"      # cr_dtimes
    $ws.Cells.Item($r, 11).Value2 = "now()"      # eff_dtimes
}

# Update the sheet selection to reflect where editing left off.
$ws.Range("D145").Select() | Out-Null
